$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new rows at the top of the data block (pushes existing rows 70-87 down to 72-89)
$ws.Rows("70:71").Insert()

# New weekly entries (date 2021-11-04 / serial 44504), mirroring the constant
# columns shared by the rest of this subset's rows.
$ws.Range("A70").Value = 2
$ws.Range("B70").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C70").Value = "Coquimbo"
$ws.Range("D70").Value = 44504
$ws.Range("E70").Value = 4
$ws.Range("F70").Value = 100112043
$ws.Range("G70").Value = "Pepino ensalada"
$ws.Range("H70").Value = "Sin especificar"
$ws.Range("I70").Value = "Primera"
$ws.Range("J70").Value = 800
$ws.Range("K70").Value = 6500
$ws.Range("L70").Value = 7000
$ws.Range("M70").Value = 6750
$ws.Range("N70").Value = "$/caja 70 unidades"
$ws.Range("O70").Value = "Provincia de Limarí"
$ws.Range("P70").Value = 96
$ws.Range("Q70").Value = 70
$ws.Range("R70").Value = "Hortaliza"

$ws.Range("A71").Value = 2
$ws.Range("B71").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C71").Value = "Coquimbo"
$ws.Range("D71").Value = 44504
$ws.Range("E71").Value = 4
$ws.Range("F71").Value = 100112043
$ws.Range("G71").Value = "Pepino ensalada"
$ws.Range("H71").Value = "Sin especificar"
$ws.Range("I71").Value = "Segunda"
$ws.Range("J71").Value = 700
$ws.Range("K71").Value = 4500
$ws.Range("L71").Value = 5000
$ws.Range("M71").Value = 4750
$ws.Range("N71").Value = "$/caja 100 unidades"
$ws.Range("O71").Value = "Provincia de Limarí"
$ws.Range("P71").Value = 48
$ws.Range("Q71").Value = 100
$ws.Range("R71").Value = "Hortaliza"
